# Unit12_Loop Invariant.pptx edit
# 1) NotesMaster date placeholder: "1/31/2024" -> "2/13/2024"
# 2) Slide 4 ("Assertion for loops"), TextBox 15 (shape id 16): add a yellow
#    highlight to the "// { x <= 0 }" loop-invariant comment run, matching
#    the highlight already present on the other loop-invariant comments.

$p = $ppt.ActivePresentation

# --- 1) Update the cached text of the NotesMaster date field -------------
$nm = $p.NotesMaster
$hf = $nm.HeadersFooters
$hf.DateAndTime.Text = "2/13/2024"

# --- 2) Highlight the last loop-invariant comment on slide 4 -------------
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("TextBox 15")
$tr = $sh.TextFrame.TextRange
$fullText = $tr.Text
$target = "// { x <= 0 }"
$pos = $fullText.IndexOf($target)
$run = $tr.Characters($pos + 1, $target.Length)
$run.Font.Highlight.RGB = 65535
